$wb = $excel.ActiveWorkbook

# Mapping: worksheet index (1-based) -> new B102 value, and row 103 values (A103 is fixed date serial 45966)
$updates = @(
    @{ Sheet = 1; B102 = 870350;  B103 = 849116  },
    @{ Sheet = 2; B102 = 1263904; B103 = 1333786 },
    @{ Sheet = 3; B102 = 134153;  B103 = 132363  },
    @{ Sheet = 4; B102 = 161188;  B103 = 159638  }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    # Update existing B102 (was 0) to the new remn_amt value
    $ws.Cells.Item(102, 2).Value = $u.B102

    # Add new row 103: date serial 45966 in A103 (same number format as A102), and remn_amt in B103
    $ws.Cells.Item(103, 1).Value = 45966
    $ws.Cells.Item(103, 1).NumberFormat = $ws.Cells.Item(102, 1).NumberFormat
    $ws.Cells.Item(103, 2).Value = $u.B103
}
